$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Two new LeetCode problems ("Repeated DNA Sequence" and "Longest
# Harmonious Sequence") were solved and inserted as rows 54-55.  The
# existing rows that used to occupy 54-58 move down two rows to 56-60.
# We reproduce that by first pushing the existing row contents (B:E)
# down two rows (58->60, 57->59, 56->58, 55->57, 54->56), starting
# from the bottom so nothing is overwritten before it is copied, and
# only then filling rows 54-55 with the new problem data.
# ------------------------------------------------------------------

# Shift row 58 -> row 60
$ws.Range("B60").Value = $ws.Range("B58").Value2
$ws.Range("C60").Value = $ws.Range("C58").Value2
$ws.Range("D60").Value = $ws.Range("D58").Value2
$ws.Range("E60").Value = $ws.Range("E58").Value2

# Shift row 57 -> row 59
$ws.Range("B59").Value = $ws.Range("B57").Value2
$ws.Range("C59").Value = $ws.Range("C57").Value2
$ws.Range("D59").Value = $ws.Range("D57").Value2
$ws.Range("E59").Value = $ws.Range("E57").Value2

# Shift row 56 -> row 58
$ws.Range("B58").Value = $ws.Range("B56").Value2
$ws.Range("C58").Value = $ws.Range("C56").Value2
$ws.Range("D58").Value = $ws.Range("D56").Value2
$ws.Range("E58").Value = $ws.Range("E56").Value2

# Shift row 55 -> row 57
$ws.Range("B57").Value = $ws.Range("B55").Value2
$ws.Range("C57").Value = $ws.Range("C55").Value2
$ws.Range("D57").Value = $ws.Range("D55").Value2
$ws.Range("E57").Value = $ws.Range("E55").Value2

# Shift row 54 -> row 56
$ws.Range("B56").Value = $ws.Range("B54").Value2
$ws.Range("C56").Value = $ws.Range("C54").Value2
$ws.Range("D56").Value = $ws.Range("D54").Value2
$ws.Range("E56").Value = $ws.Range("E54").Value2

# ------------------------------------------------------------------
# Row 54: "Repeated DNA Sequence" (Hashing, Medium, Done)
# ------------------------------------------------------------------
$ws.Range("C54").Value = "Repeated DNA Sequence"
$ws.Range("E54").Value = "Done"
$ws.Range("F51").Copy()
$ws.Range("F54").PasteSpecial(-4122)
$ws.Range("F54").Value = (Get-Date -Year 2025 -Month 8 -Day 27 -Hour 0 -Minute 0 -Second 0)
$ws.Range("G54").Value = "O(10 * n)"
$ws.Range("H54").Value = "O(10 * n)"
$ws.Range("I54").Value = "Set for seen & repeated"

# ------------------------------------------------------------------
# Row 55: "Longest Harmonious Sequence" (Easy, Done)
# ------------------------------------------------------------------
$ws.Range("C55").Value = "Longest Harmonious Sequence"
$ws.Range("D55").Value = "Easy"
$ws.Range("E55").Value = "Done"
$ws.Range("F52").Copy()
$ws.Range("F55").PasteSpecial(-4122)
$ws.Range("F55").Value = (Get-Date -Year 2025 -Month 8 -Day 27 -Hour 0 -Minute 0 -Second 0)
$ws.Range("G55").Value = "O(n)"
$ws.Range("H55").Value = "O(n)"
$ws.Range("I55").Value = "Using Frequency"

# ------------------------------------------------------------------
# Update the view: scrolled down a bit further and last edited cell
# was I56 (the Approach Summary of the "Count Occurrences of
# Anagrams" row, now sitting at row 57... matches the saved selection)
# ------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 41
$ws.Range("I56").Select()
